$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5913.8237
$ws.Range("I40").Value = 4294.9165
$ws.Range("J40").Value = 9799.200000000001
$ws.Range("K40").Value = 4294.9165
$ws.Range("L40").Value = 9799.200000000001
$ws.Range("M40").Value = -4119.9165
$ws.Range("N40").Value = -10149.2
$ws.Range("H55").Value = 96
$ws.Range("H68").Value = 55147.5
$ws.Range("J68").Value = 55147.5
$ws.Range("L68").Value = 55147.5
$ws.Range("N68").Value = -56645.5
$ws.Range("H71").Value = 55147.5
$ws.Range("J71").Value = 55147.5
$ws.Range("L71").Value = 165442.5
$ws.Range("N71").Value = -172930.5
$ws.Range("H74").Value = 13318.333
$ws.Range("I74").Value = 3964
$ws.Range("K74").Value = 3964
$ws.Range("M74").Value = -3028
$ws.Range("H77").Value = 13318.333
$ws.Range("I77").Value = 3964
$ws.Range("K77").Value = 19820
$ws.Range("M77").Value = -15140
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4846.9414
$ws.Range("I2").Value = 2339.2307
$ws.Range("J2").Value = 12997
$ws.Range("K2").Value = 2339.2307
$ws.Range("L2").Value = 12997
$ws.Range("M2").Value = -2226.2307
$ws.Range("N2").Value = -13223
$ws.Range("H32").Value = 4595.5713
$ws.Range("I32").Value = 4595.5713
$ws.Range("K32").Value = 4595.5713
$ws.Range("M32").Value = -4308.5713
$ws.Range("H36").Value = 8494
$ws.Range("I36").Value = 8494
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 8494
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -8148
$ws.Range("N36").ClearContents()
$ws.Range("H61").Value = 6406.4443
$ws.Range("I61").Value = 6216.6665
$ws.Range("K61").Value = 6216.6665
$ws.Range("M61").Value = -6004.6665
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 6558.4
$ws.Range("J74").Value = 8017.2
$ws.Range("L74").Value = 8017.2
$ws.Range("N74").Value = -9765.200000000001
$ws.Range("H77").Value = 6558.4
$ws.Range("J77").Value = 8017.2
$ws.Range("L77").Value = 40086
$ws.Range("N77").Value = -48822
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H88").Value = 1465.8889
$ws.Range("I88").Value = 1043.75
$ws.Range("J88").Value = 1803.6
$ws.Range("K88").Value = 1043.75
$ws.Range("L88").Value = 1803.6
$ws.Range("M88").Value = -637.75
$ws.Range("N88").Value = -2615.6
$ws.Range("H91").Value = 1465.8889
$ws.Range("I91").Value = 1043.75
$ws.Range("J91").Value = 1803.6
$ws.Range("K91").Value = 1043.75
$ws.Range("L91").Value = 1803.6
$ws.Range("M91").Value = 360.25
$ws.Range("N91").Value = -4611.6
$ws.Range("H110").Value = 6702.2
$ws.Range("I110").Value = 6377.75
$ws.Range("K110").Value = 6377.75
$ws.Range("M110").Value = -4332.75
$ws.Range("H116").Value = 4846.9414
$ws.Range("I116").Value = 2339.2307
$ws.Range("J116").Value = 12997
$ws.Range("K116").Value = 2339.2307
$ws.Range("L116").Value = 12997
$ws.Range("M116").Value = -45.23070000000007
$ws.Range("N116").Value = -17585
$ws.Range("H132").Value = 2944.3572
$ws.Range("I132").Value = 1485.7778
$ws.Range("K132").Value = 4457.3334
$ws.Range("M132").Value = -1927.3334
$ws.Range("H136").Value = 6406.4443
$ws.Range("I136").Value = 6216.6665
$ws.Range("K136").Value = 18649.9995
$ws.Range("M136").Value = -16099.9995
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4846.9414
$ws.Range("I3").Value = 2339.2307
$ws.Range("J3").Value = 12997
$ws.Range("K3").Value = 2339.2307
$ws.Range("L3").Value = 12997
$ws.Range("M3").Value = -2225.2307
$ws.Range("N3").Value = -13225
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2909
$ws.Range("I94").Value = 1200.125
$ws.Range("J94").Value = 4617.875
$ws.Range("K94").Value = 1200.125
$ws.Range("L94").Value = 4617.875
$ws.Range("M94").Value = -749.125
$ws.Range("N94").Value = -5519.875
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 362.63635
$ws.Range("I8").Value = 362.63635
$ws.Range("K8").Value = 1087.90905
$ws.Range("M8").Value = -948.90905
$ws.Range("H109").Value = 169597.67
$ws.Range("I109").Value = 202381.6
$ws.Range("K109").Value = 607144.8
$ws.Range("M109").Value = -606104.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H48").Value = 10041
$ws.Range("I48").Value = 10041
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 10041
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -9380
$ws.Range("N48").ClearContents()
$ws.Range("H53").Value = 19995
$ws.Range("I53").Value = 19995
$ws.Range("K53").Value = 19995
$ws.Range("M53").Value = -19477
$ws.Range("H68").Value = 6122.067
$ws.Range("I68").Value = 4364
$ws.Range("K68").Value = 4364
$ws.Range("M68").Value = -3615
$ws.Range("H71").Value = 6122.067
$ws.Range("I71").Value = 4364
$ws.Range("K71").Value = 21820
$ws.Range("M71").Value = -18076
$ws.Range("H93").Value = 1660.1
$ws.Range("I93").Value = 1871.8572
$ws.Range("J93").Value = 1166
$ws.Range("K93").Value = 1871.8572
$ws.Range("L93").Value = 1166
$ws.Range("M93").Value = -623.8571999999999
$ws.Range("N93").Value = -3662
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2999.6667
$ws.Range("I122").Value = 1999.5
$ws.Range("K122").Value = 5998.5
$ws.Range("M122").Value = -3548.5
